$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 2554225.24
$ws.Range("C9").Value = 368296.99
$ws.Range("D9").Value = 2922522.23
$ws.Range("E9").Value = 12.60202527184883
$ws.Range("F9").Value = 87.39797472815115
$ws.Range("G9").Value = -64.40593232995384
$ws.Range("H9").Value = -53.87422271691116
$ws.Range("I9").Value = -55.53230551453142
$ws.Range("J9").Value = 25323
$ws.Range("K9").Value = 1061
$ws.Range("L9").Value = 26384
